# Applies the edit described by the commit "Se resuelve el problema":
# fills in the "Ejecución de la Prueba" (row 13) timing and the
# "Desarrollo y correctivos" table rows 18-21 with real data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Row 13: Ejecución de la Prueba (Tiempo Estimado / Hora Inicio / Hora Fin) ---
$ws.Range("B13").Value = 0.027777777777777776
$ws.Range("C13").Value = 0.625
$ws.Range("D13").Value = 0.6458333333333334

# --- Rows 18-21: Desarrollo y correctivos table ---
# The shared-string table must grow in this exact order (EjercicioOIA,
# Main, EnvasadoraDeLatas, Secuencia) to match the target file's string
# indices, so the C-column text is written before anything else touches
# the shared-string pool.
$ws.Range("C18").Value = "EjercicioOIA"
$ws.Range("C19").Value = "Main"
$ws.Range("C21").Value = "EnvasadoraDeLatas"
$ws.Range("C20").Value = "Secuencia"

# Row 18 - EjercicioOIA
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 0.0006944444444444445
$ws.Range("H18").Value = 0.6458333333333334
$ws.Range("I18").Value = 0.6465277777777778
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0

# Row 19 - Main
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 0.001388888888888889
$ws.Range("H19").Value = 0.6465277777777778
$ws.Range("I19").Value = 0.6472222222222223
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0

# Row 20 - Secuencia
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 0.001388888888888889
$ws.Range("H20").Value = 0.6472222222222223
$ws.Range("I20").Value = 0.6493055555555556
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0

# Row 21 - EnvasadoraDeLatas
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 0.034722222222222224
$ws.Range("H21").Value = 0.6493055555555556
$ws.Range("I21").Value = 0.6770833333333334
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0

# F39 (=IF(E39="Completar",E39,IFERROR(E39/$E$43,"Error"))) is the lone
# member of a single-cell shared-formula group; touching it forces a
# fresh evaluation against the new E39/E43 inputs above.
$ws.Range("F39").Formula = $ws.Range("F39").Formula

# --- Selection / view state left where the author ended up editing ---
$ws.Range("A15").Select
$ws.Range("M18").Select

$wb.Application.Calculate()
